$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume(1h)) hold text-like values (e.g. "37.009.01",
# "1.00", "0.0₃0863") that must stay as literal text. Excel's automatic
# "smart" value parsing would otherwise coerce numeric-looking strings into
# real numbers (dropping formatting such as trailing zeros) or mis-parse the
# subscript-digit strings entirely. Force the whole data range to Text format
# before writing any new values so everything round-trips as plain text.
$ws.Range("D2:E51").NumberFormat = "@"

# subscript-three character used in row 20's price ("0.0\x{2083}0863")
$sub3 = [char]0x2083

function Set-Row($row, $price, $volume) {
    if ($price -ne $null) {
        $ws.Cells.Item($row, 4).Value = $price
    }
    if ($volume -ne $null) {
        $ws.Cells.Item($row, 5).Value = $volume
    }
}

Set-Row 2  "37.247.92"  "  +1.83%  "
Set-Row 3  "1.991.05"   "  +1.66%  "
Set-Row 4  $null         "  -0.23%  "
Set-Row 5  "245.67"     "  +0.56%  "
Set-Row 6  "0.631"      "  +2.26%  "
Set-Row 7  "61.71"      "  +5.29%  "
Set-Row 8  "0.999"      "  -0.15%  "
Set-Row 9  "0.384"      "  +1.26%  "
Set-Row 10 "0.0802"     "  -0.79%  "
Set-Row 11 "0.104"      "  +0.22%  "
Set-Row 12 "14.91"      "  +8.61%  "
Set-Row 13 "22.45"      "  +1.45%  "
Set-Row 14 "0.846"      "  +1.79%  "
Set-Row 15 "2.271.44"   "  +1.14%  "
Set-Row 16 $null         "  +2.69%  "
Set-Row 17 "1.990.83"   "  +1.73%  "
Set-Row 18 "37.096.03"  "  +1.65%  "
Set-Row 19 "70.26"      "  +0.74%  "
Set-Row 20 "0.0$($sub3)0863" "  +0.94%  "
Set-Row 21 "5.18"       "  +2.44%  "
Set-Row 22 "230.60"     "  +0.91%  "
Set-Row 23 $null         "  +0.19%  "
Set-Row 24 $null         "  +2.46%  "
Set-Row 25 "2.37"       "  +1.01%  "
Set-Row 26 $null         "  +4.49%  "
Set-Row 27 "9.32"       "  +0.59%  "
Set-Row 28 "163.75"     "  +2.11%  "
Set-Row 29 "19.68"      "  +1.13%  "
Set-Row 30 "1.37"       "  +18.20%  "
Set-Row 31 "0.122"      "  +1.65%  "
Set-Row 32 "4.86"       "  +3.09%  "
Set-Row 33 "0.0625"     "  +0.88%  "
Set-Row 34 "4.60"       "  +6.35%  "
Set-Row 35 "2.31"       "  +3.06%  "
Set-Row 36 "0.999"      "  -0.34%  "

# Rows 37/38 swap coin identity (WEMIXToken <-> RenderToken) in addition to
# getting refreshed price/volume figures. Rank numbers in column A stay put.
$ws.Cells.Item(37, 2).Value = "RenderToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(37, 4).Value = "3.36"
$ws.Cells.Item(37, 5).Value = "  -0.72%  "

$ws.Cells.Item(38, 2).Value = "WEMIXToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(38, 4).Value = "1.79"
$ws.Cells.Item(38, 5).Value = "  +0.87%  "

Set-Row 39 "5.50"       "  -3.33%  "
Set-Row 40 "0.0980"     "  +0.09%  "
Set-Row 41 $null         "  +1.42%  "
Set-Row 42 "1.18"       "  +0.89%  "
Set-Row 43 "0.0214"     "  +0.99%  "
Set-Row 44 "16.72"      "  +4.26%  "
Set-Row 45 "1.377.03"   "  +0.65%  "
Set-Row 46 "90.44"      "  +2.83%  "
Set-Row 47 "1.04"       "  +0.70%  "
Set-Row 48 "7.23"       "  +0.92%  "
Set-Row 49 "2.83"       "  +0.16%  "
Set-Row 50 "46.64"      "  +6.79%  "
Set-Row 51 "1.99"       "  +11.67%  "
